$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update component unit price values for rows 35 and 36 (column C)
$ws.Range("C35").Value = 2.9157999999999999
$ws.Range("C36").Value = 3.2934000000000001

# Update the active selection to match the saved cursor position
$ws.Range("H32").Select()
